# "Updated notebook, reran simulation"
#
# The underlying Jupyter notebook's method list grew from 28 to 30 entries:
# two new methods, "Holden" and "Rizzie Spiral", were inserted right after
# "Spiral5", and "Thomas Hex" was renamed to "Matthies Hex". Re-running the
# simulation regenerated the whole confusion-matrix sheet against the new,
# longer method list, so every method originally at list position >= 2
# shifted down two rows, and the table grew two new rows (30/31, index
# 28/29) at the bottom to cover the two extra methods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated method list (row labels), in their new order/positions.
$methods = @(
    "HKL",
    "Spiral5",
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

# Rows 2..31 hold one method each (A = index 0..29, B = method name).
# Row 2 (index 0, "HKL") and row 3 (index 1, "Spiral5") are already
# correct and unchanged; refresh B4:B31 for indices 2..29, and extend
# the table with the two new rows 30/31 (A/B plus the C:T data columns).
for ($i = 2; $i -le 29; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $methods[$i]
}

# New rows 30 and 31: index column, method name, and a full row of 1s
# across the data columns C:T (columns 3..20), matching every other
# data row in the sheet.
for ($i = 28; $i -le 29; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $methods[$i]
    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }

    # Match the bold/bordered/centered style already used on every other
    # index-column (A) cell in the sheet.
    $idxCell = $ws.Cells.Item($row, 1)
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1
}
